# Actualización automática 2025-08-27 10:00:10
# Updates the "CUMPLIMIENTO MENSUAL" sheet with refreshed VENTA figures
# and the derived POR CUMPLIR / CUMPLIMIENTO totals.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("CUMPLIMIENTO MENSUAL")

# Row 2 - GRUPO: OTROS
$ws.Range("D2").Value = 72466.41
$ws.Range("E2").Value = -72466.41

# Row 3 - GRUPO: PORCELANATO
$ws.Range("D3").Value = 3287.87
$ws.Range("E3").Value = 12183.6893
$ws.Range("F3").Value = 0.2125105773921572

# Row 4 - TOTAL
$ws.Range("D4").Value = 75754.28
$ws.Range("E4").Value = -60282.72070000001
$ws.Range("F4").Value = 4.896357149986814
